$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.253.93"
$ws.Range("E2").Value = "  -1.63%  "
$ws.Range("D3").Value = "2.248.51"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.86"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.49"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -4.63%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.621"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.28"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.53%  "
$ws.Range("E11").Value = "  -2.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.21"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("E13").Value = "  -2.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.59"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.45%  "
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").Value = "2.255.72"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").Value = "42.159.59"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "0.0₃0988"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.49"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.15"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.23"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.86"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.84"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +36.99%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.62"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.79%  "
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.24"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.41"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0825"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.125"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.11"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.78%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.125"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.25"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +10.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.49"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0316"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.11"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("E39").Value = "  -3.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.79"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "62.25"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "107.27"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -5.02%  "
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.71"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.27%  "
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("E47").Value = "  -3.28%  "
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("E49").Value = "  +2.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.18"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -8.70%  "
$ws.Range("E51").Value = "  -3.48%  "
